$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest scrape: row number => @(new Price text, new Volume(1h) text)
$updates = @{
    2 = @('27.735.59', '  -4.56%  ')
    3 = @('1.720.34', '  -5.29%  ')
    4 = @('1.002', '  -0.35%  ')
    5 = @('223.90', '  -3.82%  ')
    6 = @('0.5632', '  -4.00%  ')
    7 = @('1.003', '  -0.32%  ')
    8 = @('0.2696', '  -0.71%  ')
    9 = @('22.63', '  -0.65%  ')
    10 = @('0.06534', '  -3.43%  ')
    11 = @('0.07502', '  -0.32%  ')
    12 = @('1.720.13', '  -5.25%  ')
    13 = @('4.649', '  +0.50%  ')
    14 = @('0.5916', '  -4.18%  ')
    15 = @('1.956.41', '  -4.39%  ')
    16 = @('73.55', '  -1.35%  ')
    17 = @('0.000008479', '  -9.93%  ')
    18 = @('27.718.72', '  -3.96%  ')
    19 = @('5.233', '  -3.24%  ')
    20 = @('1.004', '  -0.17%  ')
    21 = @('11.15', '  -2.01%  ')
    22 = @('198.82', '  -3.84%  ')
    23 = @('6.485', '  -3.65%  ')
    24 = @('1.004', '  -0.32%  ')
    25 = @('149.11', '  -3.22%  ')
    26 = @('7.913', '  +1.96%  ')
    27 = @('0.1207', '  -3.82%  ')
    28 = @('15.98', '  -1.17%  ')
    29 = @('1.363', '  -3.11%  ')
    30 = @('0.06073', '  -3.51%  ')
    31 = @('1.379', '  -3.32%  ')
    32 = @('3.674', '  -0.46%  ')
    33 = @('3.662', '  -0.44%  ')
    34 = @('1.663', '  -1.12%  ')
    35 = @('1.022', '  -2.37%  ')
    36 = @('0.6397', '  +1.03%  ')
    37 = @('2.421', '  -4.36%  ')
    38 = @('2.675', '  -2.40%  ')
    39 = @('0.01659', '  -2.43%  ')
    40 = @('1.108.19', '  -1.89%  ')
    41 = @('6.113', '  -4.01%  ')
    42 = @('0.8731', '  +1.68%  ')
    43 = @('1.004', '  -0.22%  ')
    44 = @('99.00', '  -0.97%  ')
    45 = @('1.869.23', '  -5.05%  ')
    46 = @('58.56', '  -2.44%  ')
    47 = @('0.00000000106', '  -5.79%  ')
    48 = @('1.536', '  -1.61%  ')
    49 = @('8.187', '  -0.05%  ')
    50 = @('0.05356', '  -2.26%  ')
    51 = @('0.4405', '  -2.96%  ')
}

foreach ($row in $updates.Keys) {
    $priceText = $updates[$row][0]
    $volText = $updates[$row][1]

    $priceCell = $ws.Range("D$row")
    if ($priceText -match '^-?\d+(\.\d+)?$') {
        # The source values are plain text (e.g. "223.90", keeping the
        # trailing zero) rather than numbers, so force text entry the same
        # way a quote-prefixed literal does in Excel, then drop back to the
        # default style so no extra number formatting sticks around.
        $priceCell.Value = "'" + $priceText
        $priceCell.Style = "Normal"
    } else {
        $priceCell.Value = $priceText
    }

    $ws.Range("E$row").Value = $volText
}
